# The deck swapped the roles of its two embedded themes: the slide
# master ("theme1.xml") used to carry the "Integral" color scheme and
# now carries the stock "Office Theme" color scheme (while the
# secondary theme keeps the "Integral" colors). The PowerPoint object
# model exposes the applied theme's colors on the slide master via
# SlideMaster.Theme.ThemeColorScheme, so re-point every theme color
# slot (except Dark1/Light1, which are unchanged black/white) from the
# "Integral" values to the standard "Office" values.

function ConvertTo-VbaBgr([string]$hex) {
    # VBA's ColorFormat.RGB packs colors as 0x00BBGGRR (reverse byte
    # order of the usual RRGGBB hex notation used in OOXML).
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# index -> target "Office" theme color (RRGGBB)
$officeColors = @{
    3  = "44546A"  # Dark 2
    4  = "E7E6E6"  # Light 2
    5  = "5B9BD5"  # Accent 1
    6  = "ED7D31"  # Accent 2
    7  = "A5A5A5"  # Accent 3
    8  = "FFC000"  # Accent 4
    9  = "4472C4"  # Accent 5
    10 = "70AD47"  # Accent 6
    11 = "0563C1"  # Hyperlink
    12 = "954F72"  # Followed Hyperlink
}

foreach ($idx in $officeColors.Keys | Sort-Object) {
    $tcs.Colors($idx).RGB = ConvertTo-VbaBgr $officeColors[$idx]
}
